$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add the new data row (row 19) - extends the used range from A1:G18 to A1:G19
$ws.Range("A19").Value = 2
$ws.Range("B19").Value = 1
$ws.Range("C19").Value = 1
$ws.Range("D19").Value = 0.00000079200000000
$ws.Range("E19").Value = 0.00000000038500000
$ws.Range("F19").Value = 0.00649000000000000
$ws.Range("G19").Value = 0.00000202000000000

# Match the scientific-notation number format used by the other "p-value" style
# columns (D, E, G) further up the sheet (style index 1 -> 0.00E+00)
$ws.Range("D19").NumberFormat = "0.00E+00"
$ws.Range("E19").NumberFormat = "0.00E+00"
$ws.Range("G19").NumberFormat = "0.00E+00"

# Update the view: scroll so row 4 is the top visible row, and move the
# active selection to the newly added cell E19
$ws.Application.Goto($ws.Range("A4"))
$ws.Range("E19").Select()
